$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 200154
$ws.Range("I9").Value = 333456.66
$ws.Range("K9").Value = 333456.66
$ws.Range("M9").Value = -333287.66
$ws.Range("H17").Value = 1831.2954
$ws.Range("J17").Value = 1831.2954
$ws.Range("L17").Value = 5493.8862
$ws.Range("N17").Value = -5829.8862
$ws.Range("H32").Value = 900
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -1652
$ws.Range("H41").Value = 510.25
$ws.Range("I41").Value = 495.5
$ws.Range("J41").Value = 525
$ws.Range("K41").Value = 495.5
$ws.Range("L41").Value = 525
$ws.Range("M41").Value = -55.5
$ws.Range("N41").Value = -1405
$ws.Range("H64").Value = 8000
$ws.Range("I64").Value = 8000
$ws.Range("K64").Value = 8000
$ws.Range("M64").Value = -7752
$ws.Range("H67").Value = 8000
$ws.Range("I67").Value = 8000
$ws.Range("K67").Value = 8000
$ws.Range("M67").Value = -7142
$ws.Range("H74").Value = 10333.333
$ws.Range("I74").Value = 7500
$ws.Range("K74").Value = 7500
$ws.Range("M74").Value = -6564
$ws.Range("H77").Value = 10333.333
$ws.Range("I77").Value = 7500
$ws.Range("K77").Value = 37500
$ws.Range("M77").Value = -32820
$ws.Range("H125").Value = 474.5
$ws.Range("I125").Value = 474.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 4270.5
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -1810.5
$ws.Range("H137").Value = 2260.5386
$ws.Range("I137").Value = 2376.0908
$ws.Range("K137").Value = 7128.2724
$ws.Range("M137").Value = -4578.2724
$ws.Range("H138").Value = 4467.143
$ws.Range("I138").Value = 3516.2
$ws.Range("J138").Value = 4995.4443
$ws.Range("K138").Value = 10548.6
$ws.Range("L138").Value = 14986.3329
$ws.Range("M138").Value = -5408.599999999999
$ws.Range("N138").Value = -25266.3329

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5940.25
$ws.Range("J2").Value = 5304.3335
$ws.Range("L2").Value = 5304.3335
$ws.Range("N2").Value = -5530.3335
$ws.Range("H61").Value = 4642.75
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 4223.189
$ws.Range("I74").Value = 4705.1562
$ws.Range("J74").Value = 1138.6
$ws.Range("K74").Value = 4705.1562
$ws.Range("L74").Value = 1138.6
$ws.Range("M74").Value = -3831.1562
$ws.Range("N74").Value = -2886.6
$ws.Range("H77").Value = 4223.189
$ws.Range("I77").Value = 4705.1562
$ws.Range("J77").Value = 1138.6
$ws.Range("K77").Value = 23525.781
$ws.Range("L77").Value = 5693
$ws.Range("M77").Value = -19157.781
$ws.Range("N77").Value = -14429
$ws.Range("H116").Value = 5940.25
$ws.Range("J116").Value = 5304.3335
$ws.Range("L116").Value = 5304.3335
$ws.Range("N116").Value = -9892.333500000001
$ws.Range("H136").Value = 4642.75
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5940.25
$ws.Range("J3").Value = 5304.3335
$ws.Range("L3").Value = 5304.3335
$ws.Range("N3").Value = -5532.3335
$ws.Range("H134").Value = 7028.6
$ws.Range("I134").Value = 2732.5
$ws.Range("J134").Value = 13472.75
$ws.Range("K134").Value = 8197.5
$ws.Range("L134").Value = 40418.25
$ws.Range("M134").Value = -5662.5
$ws.Range("N134").Value = -45488.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1690
$ws.Range("I31").Value = 1935
$ws.Range("J31").Value = 1200
$ws.Range("K31").Value = 1935
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = -1640
$ws.Range("N31").Value = -1790
$ws.Range("H34").Value = 1690
$ws.Range("I34").Value = 1935
$ws.Range("J34").Value = 1200
$ws.Range("K34").Value = 1935
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = -1733
$ws.Range("N34").Value = -1604
$ws.Range("H58").Value = 3299.6296
$ws.Range("I58").Value = 3386.3635
$ws.Range("K58").Value = 3386.3635
$ws.Range("M58").Value = -3183.3635
$ws.Range("H99").Value = 2707
$ws.Range("I99").Value = 2707
$ws.Range("K99").Value = 2707
$ws.Range("M99").Value = -1209
$ws.Range("H126").Value = 2707
$ws.Range("I126").Value = 2707
$ws.Range("K126").Value = 8121
$ws.Range("M126").Value = -5651
$ws.Range("H136").Value = 3299.6296
$ws.Range("I136").Value = 3386.3635
$ws.Range("K136").Value = 10159.0905
$ws.Range("M136").Value = -7609.0905

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2133.45
$ws.Range("I131").Value = 2112
$ws.Range("J131").Value = 2138
$ws.Range("K131").Value = 6336
$ws.Range("L131").Value = 6414
$ws.Range("M131").Value = -1296
$ws.Range("N131").Value = -16494
$ws.Range("H134").Value = 3461.375
$ws.Range("I134").Value = 3482.1667
$ws.Range("K134").Value = 10446.5001
$ws.Range("M134").Value = -5376.500100000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4429.143
$ws.Range("I122").Value = 3900.8
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 11702.4
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -9252.400000000001
$ws.Range("N122").Value = -22150
$ws.Range("H126").Value = 4499.857
$ws.Range("I126").Value = 4499.857
$ws.Range("K126").Value = 13499.571
$ws.Range("M126").Value = -11029.571

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H122").Value = 3217.6667
$ws.Range("J122").Value = 2250
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 44990.332
$ws.Range("J74").Value = 44986
$ws.Range("L74").Value = 44986
$ws.Range("N74").Value = -46858
$ws.Range("H77").Value = 44990.332
$ws.Range("J77").Value = 44986
$ws.Range("L77").Value = 134958
$ws.Range("N77").Value = -144318
$ws.Range("H126").Value = 3400
$ws.Range("I126").Value = 3750
$ws.Range("J126").Value = 3166.6667
$ws.Range("K126").Value = 11250
$ws.Range("L126").Value = 9500.000100000001
$ws.Range("M126").Value = -8780
$ws.Range("N126").Value = -14440.0001
